$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Resize the three inline pictures (points == EMU / 12700)
# ------------------------------------------------------------------
# cotorra.jpg : 4572000 x 2571750 EMU -> 2540000 x 1428750 EMU
$shp = $d.InlineShapes.Item(1)
$shp.Width  = 200.0
$shp.Height = 112.5

# vinagrera.jpg : 5334000 x 3211512 EMU -> 2373010 x 1428750 EMU
$shp = $d.InlineShapes.Item(2)
$shp.Width  = 186.85118110236
$shp.Height = 112.5

# unnamed-chunk-2-1.png : 4620126 x 3696101 EMU -> 3810000 x 3048000 EMU
# (this is the 4th InlineShape in the document; the 3rd -
#  unnamed-chunk-1-1.png - keeps its original size)
$shp = $d.InlineShapes.Item(4)
$shp.Width  = 300.0
$shp.Height = 240.0

# ------------------------------------------------------------------
# 2) Wrap the scientific names in parentheses
# ------------------------------------------------------------------
# "Cotorra argentina Myiopsitta monachus" ->
# "Cotorra argentina (Myiopsitta monachus)"
$r = $d.Content
$null = $r.Find.Execute("Cotorra argentina ", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Cotorra argentina (", 2)

$r = $d.Content
$null = $r.Find.Execute("Myiopsitta monachus", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter(")")

# "Vinagrera Oxalis pes-caprae" -> "Vinagrera (Oxalis pes-caprae)"
$r = $d.Content
$null = $r.Find.Execute("Vinagrera ", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Vinagrera (", 2)

$r = $d.Content
$null = $r.Find.Execute("Oxalis pes-caprae", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertAfter(")")
